$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rules")

# Update the greeting text for rule R10 (cell E8) and select that cell,
# mirroring the edit made directly in the workbook.
$ws.Range("E8").Value = "GIT UPDATE"
$ws.Range("E8").Select()
